# Actualización automática 2025-09-12 15:40:09
#
# A new salesperson/client, "BRITO MORALES MARIA SOLEDAD", needs to be
# inserted into the roster right before "CARAVEDO PAZMIÑO  JAHAIRA PAMELA"
# (i.e. as the new row 13) on both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets. All existing rows from 13 onward shift down by
# one, the new row gets zeroed-out totals, and the running "X de 41"
# tally row (now containing 42 people instead of 41) is updated to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (columns A:R, data rows 2:42, tally row 43)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a blank row at position 13, pushing CARAVEDO.. and everybody
# below it down by one row.
$ws1.Rows.Item(13).Insert()

$ws1.Range("A13").Value = "OFICINA-CATAECSA"
$ws1.Range("B13").Value = "BRITO MORALES MARIA SOLEDAD"
$ws1.Range("C13:R13").Value = 0

# The tally row (now row 44) reads like "0 de 41" / "3 de 41" etc. Bump
# the denominator to 42 now that there is one more person in the list.
foreach ($col in @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")) {
    $cell = $ws1.Range($col + "44")
    $cell.Value = ($cell.Text -replace "de 41", "de 42")
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (columns A:G, data rows 2:42, total row 43)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(13).Insert()

$ws2.Range("A13").Value = "OFICINA-CATAECSA"
$ws2.Range("B13").Value = "BRITO MORALES MARIA SOLEDAD"
$ws2.Range("C13:G13").Value = 0

# Keep the originally-active sheet selected (avoid leaving a different
# tab "active" just because we touched it last).
$ws1.Select()
